$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old column C ("EachCalculateProfit") is dropped entirely; the old
# column D ("TodayClose") slides left into C, taking its header+style with it.
$ws.Columns("C").Delete()

# Refresh the Stock_Id / TodayClose data (rows 2-10, A = running index,
# B = Stock_Id, C = TodayClose). A few TodayClose values are "xx.xx"-style
# text rather than numbers - enter those with a leading apostrophe so Excel
# keeps them as text instead of coercing to a plain number.
$row = 2
$ws.Cells.Item($row, 1).Value = 0
$ws.Cells.Item($row, 2).Value = 1711
$ws.Cells.Item($row, 3).Value = 28
$row++

$ws.Cells.Item($row, 1).Value = 1
$ws.Cells.Item($row, 2).Value = 2436
$ws.Cells.Item($row, 3).Value = 102.5
$row++

$ws.Cells.Item($row, 1).Value = 2
$ws.Cells.Item($row, 2).Value = 3033
$ws.Cells.Item($row, 3).Value = 31.65
$row++

$ws.Cells.Item($row, 1).Value = 3
$ws.Cells.Item($row, 2).Value = 3035
$ws.Cells.Item($row, 3).Value = 189
$row++

$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = 3141
$ws.Cells.Item($row, 3).Value = "'267.00"
$row++

$ws.Cells.Item($row, 1).Value = 5
$ws.Cells.Item($row, 2).Value = 3189
$ws.Cells.Item($row, 3).Value = 251
$row++

# Rows 8-10 are brand new - extend column A's bordered/bold header-row style
# down into them by copying the format from the row above before writing.
$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = 3588
$ws.Cells.Item($row, 3).Value = 161
$row++

$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 9
$ws.Cells.Item($row, 2).Value = 6104
$ws.Cells.Item($row, 3).Value = "'165.50"
$row++

$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = 6411
$ws.Cells.Item($row, 3).Value = "'264.00"
$row++

$excel.CutCopyMode = $false
